# Replace outliers in column B (Prices Predicted) with median values.
# Rows 2-55 on Sheet1 get updated B-column values; column A is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    64787,42743,77907,13983,81991,47714,31661,89848,62740,75382,
    50004,17165,64556,43624,57825,67518,43915,58663,43029,89404,
    57020,63293,55776,58211,32394,54721,30805,61866,32072,86546,
    55407,106390,27099,56553,32889,63980,33508,59255,36459,40079,
    57364,35235,13983,86737,42959,53227,70456,33216,46507,47763,
    81834,58343,66502,39716
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}
